$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "sd"
$ws.Range("J5").Value = "Statement-non-opinion"
$ws.Range("I8").Value = "sv"
$ws.Range("J8").Value = "Statement-opinion"
$ws.Range("I9").Value = "sd"
$ws.Range("J9").Value = "Statement-non-opinion"
$ws.Range("I11").Value = "sd"
$ws.Range("J11").Value = "Statement-non-opinion"
$ws.Range("I31").Value = "aa"
$ws.Range("J31").Value = "Agree/Accept"
$ws.Range("I32").Value = "sv"
$ws.Range("J32").Value = "Statement-opinion"
$ws.Range("I35").Value = "sd"
$ws.Range("J35").Value = "Statement-non-opinion"
$ws.Range("I47").Value = "aa"
$ws.Range("J47").Value = "Agree/Accept"
$ws.Range("I54").Value = "sv"
$ws.Range("J54").Value = "Statement-opinion"
$ws.Range("I71").Value = "%"
$ws.Range("J71").Value = "Uninterpretable"
$ws.Range("I72").Value = "%"
$ws.Range("J72").Value = "Uninterpretable"
$ws.Range("I74").Value = "ba"
$ws.Range("J74").Value = "Appreciation"
$ws.Range("I82").Value = "sd"
$ws.Range("J82").Value = "Statement-non-opinion"
$ws.Range("I109").Value = "sd"
$ws.Range("J109").Value = "Statement-non-opinion"
$ws.Range("I110").Value = "sd"
$ws.Range("J110").Value = "Statement-non-opinion"
$ws.Range("I117").Value = "b"
$ws.Range("J117").Value = "Acknowledge (Backchannel)"
$ws.Range("I118").Value = "sd"
$ws.Range("J118").Value = "Statement-non-opinion"
$ws.Range("I126").Value = "aa"
$ws.Range("J126").Value = "Agree/Accept"
$ws.Range("I136").Value = "ba"
$ws.Range("J136").Value = "Appreciation"
$ws.Range("I143").Value = "sv"
$ws.Range("J143").Value = "Statement-opinion"
$ws.Range("I144").Value = "sd"
$ws.Range("J144").Value = "Statement-non-opinion"
$ws.Range("I146").Value = "aa"
$ws.Range("J146").Value = "Agree/Accept"
$ws.Range("I147").Value = "sd"
$ws.Range("J147").Value = "Statement-non-opinion"
$ws.Range("I151").Value = "aa"
$ws.Range("J151").Value = "Agree/Accept"
$ws.Range("I156").Value = "sv"
$ws.Range("J156").Value = "Statement-opinion"
$ws.Range("I162").Value = "sv"
$ws.Range("J162").Value = "Statement-opinion"
$ws.Range("I163").Value = "%"
$ws.Range("J163").Value = "Uninterpretable"
$ws.Range("I165").Value = "sv"
$ws.Range("J165").Value = "Statement-opinion"
$ws.Range("I167").Value = "sv"
$ws.Range("J167").Value = "Statement-opinion"
$ws.Range("I177").Value = "aa"
$ws.Range("J177").Value = "Agree/Accept"
$ws.Range("I180").Value = "sd"
$ws.Range("J180").Value = "Statement-non-opinion"
$ws.Range("I181").Value = "sd"
$ws.Range("J181").Value = "Statement-non-opinion"
$ws.Range("I186").Value = "aa"
$ws.Range("J186").Value = "Agree/Accept"
$ws.Range("I187").Value = "aa"
$ws.Range("J187").Value = "Agree/Accept"
$ws.Range("I197").Value = "sd"
$ws.Range("J197").Value = "Statement-non-opinion"
$ws.Range("I199").Value = "sd"
$ws.Range("J199").Value = "Statement-non-opinion"
$ws.Range("I200").Value = "sd"
$ws.Range("J200").Value = "Statement-non-opinion"
$ws.Range("I201").Value = "b"
$ws.Range("J201").Value = "Acknowledge (Backchannel)"
$ws.Range("I204").Value = "sd"
$ws.Range("J204").Value = "Statement-non-opinion"
$ws.Range("I217").Value = "sv"
$ws.Range("J217").Value = "Statement-opinion"
$ws.Range("I218").Value = "b"
$ws.Range("J218").Value = "Acknowledge (Backchannel)"
$ws.Range("I221").Value = "ba"
$ws.Range("J221").Value = "Appreciation"
$ws.Range("I235").Value = "qy"
$ws.Range("J235").Value = "Yes-No-Question"
